$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI TPM metrics (ligand/receptor expression + derived specificity/edge weights).
# Row data below: only the columns that changed per row are listed.
$updates = @(
    @{row=2; G=0.178715; H=0.536145; I=0.09904930989061336; J=0.09904930989061336; M=198.3395793333333; N=595.018738; O=0.743678971552218; P=0.7436789715522182; Q=35.44625792055666; R=319.01632128501; S=0.07366088891240828; T=0.0736608889124083}
    @{row=3; G=0.178715; H=0.536145; I=0.09904930989061336; J=0.09904930989061336; M=44.55288066666666; O=0.1670520860498112; P=0.1670520860498113; Q=7.962268068343332; R=71.66041261509; S=0.01654639383902116; T=0.01654639383902116}
    @{row=4; G=0.178715; H=0.536145; I=0.09904930989061336; J=0.09904930989061336; M=4.884127333333333; N=14.652382; O=0.01831315163817619; P=0.0183131516381762; Q=0.8728668163766665; R=7.855801347389999; S=0.001813905031683507; T=0.001813905031683508}
    @{row=5; G=0.178715; H=0.536145; I=0.09904930989061336; J=0.09904930989061336; M=6.952303666666666; N=20.856911; O=0.02606782800550416; P=0.02606782800550417; Q=1.242480949788333; R=11.182328548095; S=0.002582000374292391; T=0.002582000374292392}
    @{row=6; G=0.178715; H=0.536145; I=0.09904930989061336; J=0.09904930989061336; M=2.587242666666667; N=7.761728; O=0.009700927933647787; P=0.009700927933647789; Q=0.4623790731733333; R=4.16141165856; S=0.0009608702171263872; T=0.0009608702171263875}
    @{row=7; G=0.178715; H=0.536145; I=0.09904930989061336; J=0.09904930989061336; M=9.384400999999999; N=28.153203; O=0.0351870348206426; P=0.03518703482064261; Q=1.677133224715; R=15.094199022435; S=0.003485251516081632; T=0.003485251516081633}
    @{row=8; E=3; F=1; G=1.450498333333333; H=4.351495; I=0.8039104659046613; J=0.8039104659046612; M=198.3395793333333; N=595.018738; O=0.743678971552218; P=0.7436789715522182; Q=287.6912292570344; R=2589.22106331331; S=0.5978513085040429; T=0.597851308504043}
    @{row=9; E=3; F=1; G=1.450498333333333; H=4.351495; I=0.8039104659046613; J=0.8039104659046612; M=44.55288066666666; O=0.1670520860498112; P=0.1670520860498113; Q=64.62387915219887; R=581.61491236979; S=0.1342949203266493; T=0.1342949203266493}
    @{row=10; E=3; F=1; G=1.450498333333333; H=4.351495; I=0.8039104659046613; J=0.8039104659046612; M=4.884127333333333; N=14.652382; O=0.01831315163817619; P=0.0183131516381762; Q=7.084418556787776; R=63.75976701109; S=0.01472213426562893; T=0.01472213426562894}
    @{row=11; E=3; F=1; G=1.450498333333333; H=4.351495; I=0.8039104659046613; J=0.8039104659046612; M=6.952303666666666; N=20.856911; O=0.02606782800550416; P=0.02606782800550417; Q=10.08430488132722; R=90.75874393194501; S=0.02095619975702743; T=0.02095619975702743}
    @{row=12; E=3; F=1; G=1.450498333333333; H=4.351495; I=0.8039104659046613; J=0.8039104659046612; M=2.587242666666667; N=7.761728; O=0.009700927933647787; P=0.009700927933647789; Q=3.752791175928889; R=33.77512058336; S=0.007798677494846336; T=0.007798677494846336}
    @{row=13; E=3; F=1; G=1.450498333333333; H=4.351495; I=0.8039104659046613; J=0.8039104659046612; M=9.384400999999999; N=28.153203; O=0.0351870348206426; P=0.03518703482064261; Q=13.61205800983166; R=122.508522088485; S=0.02828722555646633; T=0.02828722555646634}
    @{row=14; E=1; F=0.3333333333333333; G=0.17509; H=0.52527; I=0.09704022420472538; J=0.09704022420472537; M=198.3395793333333; N=595.018738; O=0.743678971552218; P=0.7436789715522182; Q=34.72727694547333; R=312.54549250926; S=0.07216677413576682; T=0.07216677413576683}
    @{row=15; E=1; F=0.3333333333333333; G=0.17509; H=0.52527; I=0.09704022420472538; J=0.09704022420472537; M=44.55288066666666; O=0.1670520860498112; P=0.1670520860498113; Q=7.800763875926665; R=70.20687488333999; S=0.01621077188414076; T=0.01621077188414076}
    @{row=16; E=1; F=0.3333333333333333; G=0.17509; H=0.52527; I=0.09704022420472538; J=0.09704022420472537; M=4.884127333333333; N=14.652382; O=0.01831315163817619; P=0.0183131516381762; Q=0.8551618547933332; R=7.69645669314; S=0.001777112340863752; T=0.001777112340863752}
    @{row=17; E=1; F=0.3333333333333333; G=0.17509; H=0.52527; I=0.09704022420472538; J=0.09704022420472537; M=6.952303666666666; N=20.856911; O=0.02606782800550416; P=0.02606782800550417; Q=1.217278848996667; R=10.95550964097; S=0.002529627874184343; T=0.002529627874184343}
    @{row=18; E=1; F=0.3333333333333333; G=0.17509; H=0.52527; I=0.09704022420472538; J=0.09704022420472537; M=2.587242666666667; N=7.761728; O=0.009700927933647787; P=0.009700927933647789; Q=0.4530003185066667; R=4.07700286656; S=0.0009413802216750646; T=0.0009413802216750646}
    @{row=19; E=1; F=0.3333333333333333; G=0.17509; H=0.52527; I=0.09704022420472538; J=0.09704022420472537; M=9.384400999999999; N=28.153203; O=0.0351870348206426; P=0.03518703482064261; Q=1.64311477109; R=14.78803293981; S=0.003414557748094637; T=0.003414557748094637}
)

foreach ($u in $updates) {
    $r = $u.row
    foreach ($col in $u.Keys) {
        if ($col -eq "row") { continue }
        $ws.Range("$col$r").Value = $u[$col]
    }
}